$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates the per-observation data (species id/name/coords, etc.)
# among rows 3, 4, 5 and 7 while row 6 and all shared site/metadata columns
# stay put: new row3 <- old row4, new row4 <- old row5, new row5 <- old row7,
# new row7 <- old row3.

# --- Row 3  (<- old Row 4 content) ---
$ws.Range("A3").Value = 111741082
$ws.Range("B3").Value = 77515
$ws.Range("E3").Value = 6425
$ws.Range("F3").Value = "Garnlav"
$ws.Range("G3").Value = "Alectoria sarmentosa"
$ws.Range("H3").Value = "(Ach.) Ach."
$ws.Range("L3").Value = ""
$ws.Range("Q3").Value = 331468.5669229594
$ws.Range("R3").Value = 6627064.351006002

# --- Row 4  (<- old Row 5 content) ---
$ws.Range("A4").Value = 111741025
$ws.Range("B4").Value = 94134
$ws.Range("E4").Value = 53
$ws.Range("F4").Value = "Vedtrappmossa"
$ws.Range("G4").Value = "Crossocalyx hellerianus"
$ws.Range("H4").Value = "(Nees ex Lindenb.) Meyl."
# L4 did not exist before; create it as an (empty) present cell
$ws.Range("L4").Value = "x"
$ws.Range("L4").Value = ""
$ws.Range("L4").NumberFormat = "General"
$ws.Range("Q4").Value = 331437.2628167981
$ws.Range("R4").Value = 6627065.263253132

# --- Row 5  (<- old Row 7 content) ---
$ws.Range("A5").Value = 111741120
$ws.Range("B5").Value = 56398
$ws.Range("E5").Value = 100109
$ws.Range("F5").Value = "Tretåig hackspett"
$ws.Range("G5").Value = "Picoides tridactylus"
$ws.Range("H5").Value = "(Linnaeus, 1758)"
$ws.Range("J5").Value = ""
$ws.Range("M5").Value = "färska spår"
$ws.Range("AF5").Value = ""
$ws.Range("Q5").Value = 331468.5669229594
$ws.Range("R5").Value = 6627064.351006002

# --- Row 7  (<- old Row 3 content) ---
$ws.Range("A7").Value = 111741038
$ws.Range("B7").Value = 94134
$ws.Range("E7").Value = 53
$ws.Range("F7").Value = "Vedtrappmossa"
$ws.Range("G7").Value = "Crossocalyx hellerianus"
$ws.Range("H7").Value = "(Nees ex Lindenb.) Meyl."
# J7 did not exist before; create it as an (empty) present cell
$ws.Range("J7").Value = "x"
$ws.Range("J7").Value = ""
$ws.Range("J7").NumberFormat = "General"
$ws.Range("M7").Value = ""
# AF7 did not exist before; create it as an (empty) present cell
$ws.Range("AF7").Value = "x"
$ws.Range("AF7").Value = ""
$ws.Range("AF7").NumberFormat = "General"
$ws.Range("Q7").Value = 331443.3172632467
$ws.Range("R7").Value = 6627064.989183033
